# "Generate Report for Handoff"
# Update the "Latest Handoff Datetime" for the 3885a69e-...md source file
# (row 4) on each locale report sheet, reflecting a fresh handoff run.
$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("D4").Value = "2016-03-10 05:33:09"
$dede.Range("D4").Value = "2016-03-10 05:33:19"
